$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add "Percentage" and "Rank" columns after "Total Marks" (AZ) ---
# Copy the header formatting (bold, centered, bordered) from the existing AZ1
# header cell onto the two new header cells before writing their text so the
# new cells pick up the same cell style index used by the rest of row 1.
$ws.Range("AZ1").Copy() | Out-Null
$ws.Range("BA1:BB1").PasteSpecial(-4122) | Out-Null

$ws.Range("BA1").Value = "Percentage"
$ws.Range("BB1").Value = "Rank"

# --- Data rows: Percentage = Total Weightage (AX) / Total Marks (AZ) * 100 ---
# --- Rank = position of Total Weightage (AX) among all students, 1 = highest ---
$ws.Range("BA2").Value = 37.68888888888889
$ws.Range("BB2").Value = 9

$ws.Range("BA3").Value = 44.53333333333334
$ws.Range("BB3").Value = 3

$ws.Range("BA4").Value = 13.68888888888889
$ws.Range("BB4").Value = 18

$ws.Range("BA5").Value = 42.26666666666667
$ws.Range("BB5").Value = 6

$ws.Range("BA6").Value = 37.88888888888889
$ws.Range("BB6").Value = 8

$ws.Range("BA7").Value = 28.75555555555556
$ws.Range("BB7").Value = 15

$ws.Range("BA8").Value = 40
$ws.Range("BB8").Value = 7

$ws.Range("BA9").Value = 47.86666666666667
$ws.Range("BB9").Value = 2

$ws.Range("BA10").Value = 43.86666666666667
$ws.Range("BB10").Value = 4

$ws.Range("BA11").Value = 27.22222222222222
$ws.Range("BB11").Value = 17

$ws.Range("BA12").Value = 37.48888888888889
$ws.Range("BB12").Value = 10

$ws.Range("BA13").Value = 49.37777777777777
$ws.Range("BB13").Value = 1

$ws.Range("BA14").Value = 42.28888888888889
$ws.Range("BB14").Value = 5

$ws.Range("BA15").Value = 35.2888888888889
$ws.Range("BB15").Value = 11

$ws.Range("BA16").Value = 31.71333333333333
$ws.Range("BB16").Value = 14

$ws.Range("BA17").Value = 34.93333333333334
$ws.Range("BB17").Value = 12

$ws.Range("BA18").Value = 27.51333333333334
$ws.Range("BB18").Value = 16

$ws.Range("BA19").Value = 34.24666666666667
$ws.Range("BB19").Value = 13
